# Auto-generated edit script: updates the "Price" (D) and "Volume(1h)" (E) columns
# of the cryptos worksheet, and fixes the TheGraph/Stellar row ordering (rows 44-45).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D values are price strings that often look numeric (e.g. "399.85").
# Force text format on the whole column first so Excel does not silently
# reinterpret them as numbers (which would also strip formatting like "39.72" -> 39.72).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "56.411.54"
$ws.Range("E2").Value = "  +9.50%  "
$ws.Range("D3").Value = "3.235.43"
$ws.Range("E3").Value = "  +4.45%  "
$ws.Range("D5").Value = "399.85"
$ws.Range("E5").Value = "  +4.05%  "
$ws.Range("D6").Value = "111.43"
$ws.Range("E6").Value = "  +8.41%  "
$ws.Range("D7").Value = "0.558"
$ws.Range("E7").Value = "  +3.32%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("E9").Value = "  +7.51%  "
$ws.Range("D10").Value = "39.72"
$ws.Range("E10").Value = "  +7.48%  "
$ws.Range("D11").Value = "0.0904"
$ws.Range("E11").Value = "  +5.60%  "
$ws.Range("E12").Value = "  +2.16%  "
$ws.Range("D13").Value = "3.749.16"
$ws.Range("E13").Value = "  +4.70%  "
$ws.Range("D14").Value = "19.23"
$ws.Range("E14").Value = "  +3.24%  "
$ws.Range("D15").Value = "8.11"
$ws.Range("E15").Value = "  +3.47%  "
$ws.Range("E16").Value = "  +8.12%  "
$ws.Range("D17").Value = "3.239.97"
$ws.Range("E17").Value = "  +4.69%  "
$ws.Range("D18").Value = "10.71"
$ws.Range("E18").Value = "  -4.30%  "
$ws.Range("D19").Value = "56.275.16"
$ws.Range("E19").Value = "  +9.28%  "
$ws.Range("D20").Value = "3.43"
$ws.Range("E20").Value = "  +3.32%  "
$ws.Range("E21").Value = "  +7.32%  "
$ws.Range("D22").Value = "13.20"
$ws.Range("E22").Value = "  +6.97%  "
$ws.Range("D23").Value = "292.87"
$ws.Range("E23").Value = "  +10.19%  "
$ws.Range("D24").Value = "74.64"
$ws.Range("E24").Value = "  +6.74%  "
$ws.Range("D25").Value = "3.25"
$ws.Range("E25").Value = "  +4.05%  "
$ws.Range("D26").Value = "8.20"
$ws.Range("E26").Value = "  +1.15%  "
$ws.Range("D27").Value = "28.21"
$ws.Range("E27").Value = "  +4.48%  "
$ws.Range("D28").Value = "7.52"
$ws.Range("E28").Value = "  +3.43%  "
$ws.Range("E29").Value = "  +3.20%  "
$ws.Range("E30").Value = "  -0.08%  "
$ws.Range("D31").Value = "0.113"
$ws.Range("E31").Value = "  +5.44%  "
$ws.Range("D32").Value = "11.37"
$ws.Range("E32").Value = "  +10.27%  "
$ws.Range("D33").Value = "0.0496"
$ws.Range("E33").Value = "  +5.40%  "
$ws.Range("D34").Value = "36.95"
$ws.Range("E34").Value = "  +4.58%  "
$ws.Range("D35").Value = "2.12"
$ws.Range("E35").Value = "  +2.85%  "
$ws.Range("D36").Value = "51.49"
$ws.Range("E36").Value = "  +2.28%  "
$ws.Range("D37").Value = "3.58"
$ws.Range("E37").Value = "  +6.72%  "
$ws.Range("D38").Value = "3.12"
$ws.Range("E38").Value = "  +24.23%  "
$ws.Range("D39").Value = "1.00"
$ws.Range("E39").Value = "  +0.18%  "
$ws.Range("D40").Value = "136.53"
$ws.Range("E40").Value = "  +6.03%  "
$ws.Range("D41").Value = "1.93"
$ws.Range("E41").Value = "  +2.79%  "
$ws.Range("E42").Value = "  +10.27%  "
$ws.Range("D43").Value = "17.17"
$ws.Range("E43").Value = "  +3.88%  "
$ws.Range("D46").Value = "22.69"
$ws.Range("E46").Value = "  +1.52%  "
$ws.Range("D47").Value = "2.14"
$ws.Range("E47").Value = "  +42.73%  "
$ws.Range("D48").Value = "2.153.34"
$ws.Range("E48").Value = "  +4.89%  "
$ws.Range("E49").Value = "  +0.93%  "
$ws.Range("E50").Value = "  -0.67%  "
$ws.Range("D51").Value = "0.0361"
$ws.Range("E51").Value = "  +9.85%  "

# Rows 44 and 45 swapped which coin they describe (Stellar now sorts above TheGraph).
$ws.Range("B44").Value = "Stellar"
$ws.Range("C44").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D44").Value = "0.119"
$ws.Range("E44").Value = "  +3.09%  "
$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D45").Value = "0.286"
$ws.Range("E45").Value = "  -4.98%  "

